# Add "TODO task" questions to the Algorithms section (rows 21-24) of Sheet1.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Column A ("Topic") for the existing + new Algorithms rows (21-24) -
# reuses the existing "Algorithms" shared string (same as rows 15-20).
$ws.Cells.Item(21, 1).Value = "Algorithms"
$ws.Cells.Item(22, 1).Value = "Algorithms"
$ws.Cells.Item(23, 1).Value = "Algorithms"
$ws.Cells.Item(24, 1).Value = "Algorithms"

# New question in row 23, with its (long, multi-line) answer in column D.
$ws.Cells.Item(23, 3).Value = "What is general form of recurrance solution?"

$ws.Cells.Item(23, 4).Value = "aT(n-b)+f(n) when a > 0, b > 0 and f(n) = O(n^k) where k >= 0`nif a=1 then  O(n^(k+1)) or O(n*f(n))`nif a > 1 then O(n^K * a^(n/b)) or O(f(n) * a^(n/b))`nif a < 1 then O(n^k) or O(f(n))"
$ws.Cells.Item(23, 4).WrapText = $true
$ws.Rows.Item(23).RowHeight = 60

# New TODO question in row 24.
$ws.Cells.Item(24, 3).Value = "What is masters theorem for decreasing functions?"

# Leave the selection where the author left off, ready for the next row.
$ws.Range("A25").Select()
